$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Windows VM name typo: "window-vm-test" -> "windows-vm-test"
$ws.Range("B5").Value = "windows-vm-test"

$ws.Range("B5").Select()

$wb.Save()
